# Apply the Sat Sep 30 03:09:16 UTC 2023 GitHub Actions "cryptos list" refresh:
# updates the Price (D) and Volume(1h) (E) columns for every coin row, and
# swaps the Stellar / BinanceUSD rows (28 <-> 29) to match the new ranking order.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D (Price) stores numeric-looking text (e.g. "1.00", "26.921.20").
# Plain `.Value = "..."` would let Excel auto-coerce that into a real number
# (losing trailing zeros / thousands-separated formatting), so briefly switch
# those cells to Text format while writing, then restore the original (default)
# style so no stray formatting is left behind.
$priceCells = @("D2", "D3", "D5", "D10", "D11", "D12", "D13", "D16", "D17", "D19", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D33", "D34", "D35", "D38", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.921.20"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.671.67"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "214.81"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").Value = "20.29"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +1.67%  "
$ws.Range("D12").Value = "1.906.54"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "1.698.41"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "65.52"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "26.916.26"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("E18").Value = "  +3.97%  "
$ws.Range("D19").Value = "233.55"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "4.43"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "146.16"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("D26").Value = "7.13"
$ws.Range("E26").Value = "  +0.42%  "
$ws.Range("D27").Value = "15.96"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("B29").Value = "BinanceUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "1.459.10"
$ws.Range("E33").Value = "  -5.07%  "
$ws.Range("D34").Value = "3.13"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "0.901"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("E40").Value = "  +12.52%  "
$ws.Range("E41").Value = "  -4.27%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").Value = "66.31"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "90.76"
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "1.54"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("D49").Value = "0.102"
$ws.Range("E49").Value = "  +2.75%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "7.63"
$ws.Range("E51").Value = "  +0.30%  "

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
